$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 447.44446
$ws.Range("I9").Value = 243
$ws.Range("J9").Value = 703
$ws.Range("K9").Value = 243
$ws.Range("L9").Value = 703
$ws.Range("M9").Value = -74
$ws.Range("N9").Value = -1041
$ws.Range("H81").Value = 112500
$ws.Range("I81").Value = 75000
$ws.Range("J81").Value = 150000
$ws.Range("K81").Value = 75000
$ws.Range("L81").Value = 150000
$ws.Range("M81").Value = -74002
$ws.Range("N81").Value = -151996
$ws.Range("H82").Value = 5456
$ws.Range("I82").Value = 1821.25
$ws.Range("K82").Value = 5463.75
$ws.Range("M82").Value = -5057.75
$ws.Range("H84").Value = 112500
$ws.Range("I84").Value = 75000
$ws.Range("J84").Value = 150000
$ws.Range("K84").Value = 225000
$ws.Range("L84").Value = 450000
$ws.Range("M84").Value = -220008
$ws.Range("N84").Value = -459984
$ws.Range("H85").Value = 5456
$ws.Range("I85").Value = 1821.25
$ws.Range("K85").Value = 5463.75
$ws.Range("M85").Value = -4059.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17445.406
$ws.Range("I32").Value = 19849.377
$ws.Range("J32").Value = 5862.636
$ws.Range("K32").Value = 19849.377
$ws.Range("L32").Value = 5862.636
$ws.Range("M32").Value = -19562.377
$ws.Range("N32").Value = -6436.636
$ws.Range("H122").Value = 1836.1428
$ws.Range("I122").Value = 1588.8125
$ws.Range("K122").Value = 4766.4375
$ws.Range("M122").Value = -2316.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 17466
$ws.Range("I96").Value = 17466
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 17466
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -14720
$ws.Range("N96").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1900
$ws.Range("J4").Value = 5500
$ws.Range("L4").Value = 5500
$ws.Range("N4").Value = -5724
$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H76").Value = 14999.5
$ws.Range("I76").Value = 14999.5
$ws.Range("K76").Value = 14999.5
$ws.Range("M76").Value = -14684.5
$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H79").Value = 14999.5
$ws.Range("I79").Value = 14999.5
$ws.Range("K79").Value = 14999.5
$ws.Range("M79").Value = -13907.5
$ws.Range("H132").Value = 2095.889
$ws.Range("I132").Value = 2095.889
$ws.Range("K132").Value = 6287.667
$ws.Range("M132").Value = -3757.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 853.1429
$ws.Range("I29").Value = 343.25
$ws.Range("J29").Value = 1533
$ws.Range("K29").Value = 1029.75
$ws.Range("L29").Value = 4599
$ws.Range("M29").Value = -752.75
$ws.Range("N29").Value = -5153
$ws.Range("H46").Value = 952.6667
$ws.Range("I46").Value = 1142.5714
$ws.Range("K46").Value = 3427.7142
$ws.Range("M46").Value = -3336.7142
$ws.Range("H82").Value = 9713.429
$ws.Range("J82").Value = 9713.429
$ws.Range("L82").Value = 29140.287
$ws.Range("N82").Value = -29952.287
$ws.Range("H85").Value = 9713.429
$ws.Range("J85").Value = 9713.429
$ws.Range("L85").Value = 29140.287
$ws.Range("N85").Value = -31948.287
$ws.Range("H126").Value = 11950
$ws.Range("I126").Value = 11950
$ws.Range("K126").Value = 35850
$ws.Range("M126").Value = -30910
$ws.Range("H136").Value = 5664.6665
$ws.Range("I136").Value = 3497
$ws.Range("K136").Value = 10491
$ws.Range("M136").Value = -5391
$ws.Range("H139").Value = 7066.25
$ws.Range("I139").Value = 1266
$ws.Range("K139").Value = 3798
$ws.Range("M139").Value = 1342
$ws.Range("H140").Value = 3605.2942
$ws.Range("I140").Value = 3414.9333
$ws.Range("J140").Value = 5033
$ws.Range("K140").Value = 10244.7999
$ws.Range("L140").Value = 15099
$ws.Range("M140").Value = -5064.7999
$ws.Range("N140").Value = -25459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H80").Value = 6500.25
$ws.Range("I80").Value = 5438.0835
$ws.Range("J80").Value = 8093.5
$ws.Range("K80").Value = 5438.0835
$ws.Range("L80").Value = 8093.5
$ws.Range("M80").Value = -4440.0835
$ws.Range("N80").Value = -10089.5
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H82").Value = 48432.332
$ws.Range("J82").Value = 54499.5
$ws.Range("L82").Value = 54499.5
$ws.Range("N82").Value = -55265.5
$ws.Range("H83").Value = 6500.25
$ws.Range("I83").Value = 5438.0835
$ws.Range("J83").Value = 8093.5
$ws.Range("K83").Value = 27190.4175
$ws.Range("L83").Value = 40467.5
$ws.Range("M83").Value = -22198.4175
$ws.Range("N83").Value = -50451.5
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H85").Value = 48432.332
$ws.Range("J85").Value = 54499.5
$ws.Range("L85").Value = 54499.5
$ws.Range("N85").Value = -57151.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 17666.334
$ws.Range("J62").Value = 17666.334
$ws.Range("L62").Value = 17666.334
$ws.Range("N62").Value = -18914.334
$ws.Range("H65").Value = 17666.334
$ws.Range("J65").Value = 17666.334
$ws.Range("L65").Value = 52999.00199999999
$ws.Range("N65").Value = -59239.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 32999.668
$ws.Range("I40").Value = 49999
$ws.Range("K40").Value = 49999
$ws.Range("M40").Value = -49850
$ws.Range("H76").Value = 22497
$ws.Range("I76").Value = 14995
$ws.Range("J76").Value = 29999
$ws.Range("K76").Value = 14995
$ws.Range("L76").Value = 29999
$ws.Range("M76").Value = -14680
$ws.Range("N76").Value = -30629
$ws.Range("H79").Value = 22497
$ws.Range("I79").Value = 14995
$ws.Range("J79").Value = 29999
$ws.Range("K79").Value = 14995
$ws.Range("L79").Value = 29999
$ws.Range("M79").Value = -13903
$ws.Range("N79").Value = -32183
$ws.Range("H127").Value = 94994
